$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.00", "3.938.94").
# Force text format on the whole price column first so Excel does not
# silently coerce/normalize the values we are about to write, then restore
# the default style so the workbook format matches the original.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '61.680.42'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '3.397.74'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '408.06'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').Value = '126.61'
$ws.Range('E6').Value = '  -2.68%  '
$ws.Range('D7').Value = '0.615'
$ws.Range('E7').Value = '  -2.61%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.723'
$ws.Range('E9').Value = '  -2.96%  '
$ws.Range('D10').Value = '0.134'
$ws.Range('E10').Value = '  -10.90%  '
$ws.Range('D11').Value = '42.33'
$ws.Range('E11').Value = '  -1.41%  '
$ws.Range('D12').Value = '9.07'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').Value = '3.938.94'
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').Value = '0.0000209'
$ws.Range('E15').Value = '  -8.82%  '
$ws.Range('D16').Value = '20.32'
$ws.Range('E16').Value = '  -3.70%  '
$ws.Range('D17').Value = '3.398.32'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '1.07'
$ws.Range('E18').Value = '  +0.59%  '
$ws.Range('D19').Value = '12.20'
$ws.Range('E19').Value = '  -1.05%  '
$ws.Range('D20').Value = '61.750.04'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').Value = '481.47'
$ws.Range('E21').Value = '  +20.12%  '
$ws.Range('D22').Value = '89.52'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('E23').Value = '  +0.77%  '
$ws.Range('D24').Value = '13.09'
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('D25').Value = '3.25'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').Value = '9.35'
$ws.Range('E26').Value = '  +7.16%  '
$ws.Range('D27').Value = '33.03'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').Value = '7.84'
$ws.Range('E29').Value = '  +3.04%  '
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('D31').Value = '11.76'
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').Value = '0.167'
$ws.Range('E32').Value = '  -3.13%  '
$ws.Range('E33').Value = '  -5.98%  '
$ws.Range('D34').Value = '40.76'
$ws.Range('E34').Value = '  -5.29%  '
$ws.Range('E35').Value = '  -0.66%  '
$ws.Range('D36').Value = '55.58'
$ws.Range('E36').Value = '  +2.66%  '
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').Value = '0.330'
$ws.Range('E39').Value = '  +5.74%  '
$ws.Range('D40').Value = '3.00'
$ws.Range('E40').Value = '  +2.94%  '
$ws.Range('D41').Value = '147.90'
$ws.Range('E41').Value = '  +4.17%  '
$ws.Range('D42').Value = '3.31'
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('E43').Value = '  -0.69%  '
$ws.Range('E44').Value = '  +3.97%  '
$ws.Range('D45').Value = '2.55'
$ws.Range('E45').Value = '  +5.55%  '
$ws.Range('D46').Value = '4.16'
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('E47').Value = '  +16.67%  '
$ws.Range('D48').Value = '16.34'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').Value = '21.91'
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.142'
$ws.Range('E50').Value = '  +8.00%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '112.06'
$ws.Range('E51').Value = '  +13.64%  '

# Restore default (General) styling on the price column now that the
# literal text values are safely stored.
$priceRange.Style = "Normal"

Write-Host "Applied cryptos update"
